$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price (D) and volume-change (E) columns for this scrape run.
# Values are written as formulas producing the literal text, then converted
# to static values via Copy/PasteSpecial so cells stay plain text (matching
# the original inlineStr cells) instead of being auto-coerced to numbers.

$ws.Range("D2").Formula = '="27.520.16"'
$ws.Range("E2").Formula = '="  +0.33%  "'
$ws.Range("D3").Formula = '="1.742.21"'
$ws.Range("E3").Formula = '="  -0.36%  "'
$ws.Range("E4").Formula = '="  +0.00%  "'
$ws.Range("D5").Formula = '="323.06"'
$ws.Range("E5").Formula = '="  +0.22%  "'
$ws.Range("D7").Formula = '="0.4488"'
$ws.Range("E7").Formula = '="  +5.90%  "'
$ws.Range("D8").Formula = '="0.3522"'
$ws.Range("E8").Formula = '="  -2.15%  "'
$ws.Range("D9").Formula = '="0.07372"'
$ws.Range("E9").Formula = '="  -1.73%  "'
$ws.Range("D10").Formula = '="41.35"'
$ws.Range("E10").Formula = '="  -1.87%  "'
$ws.Range("D11").Formula = '="1.077"'
$ws.Range("E11").Formula = '="  -2.18%  "'
$ws.Range("E12").Formula = '="  +0.00%  "'
$ws.Range("D13").Formula = '="20.41"'
$ws.Range("E13").Formula = '="  -1.16%  "'
$ws.Range("E14").Formula = '="  -2.22%  "'
$ws.Range("D15").Formula = '="7.063"'
$ws.Range("E15").Formula = '="  -2.13%  "'
$ws.Range("D16").Formula = '="1.741.97"'
$ws.Range("E16").Formula = '="  -0.39%  "'
$ws.Range("D17").Formula = '="91.47"'
$ws.Range("E17").Formula = '="  -1.55%  "'
$ws.Range("E18").Formula = '="  -1.47%  "'
$ws.Range("D19").Formula = '="0.06360"'
$ws.Range("E19").Formula = '="  -0.04%  "'
$ws.Range("E20").Formula = '="  +0.05%  "'
$ws.Range("D21").Formula = '="16.75"'
$ws.Range("E21").Formula = '="  -1.48%  "'
$ws.Range("E22").Formula = '="  -2.79%  "'
$ws.Range("D23").Formula = '="27.555.44"'
$ws.Range("E23").Formula = '="  +0.26%  "'
$ws.Range("E24").Formula = '="  -0.73%  "'
$ws.Range("E25").Formula = '="  +0.55%  "'
$ws.Range("D26").Formula = '="162.04"'
$ws.Range("E26").Formula = '="  -0.02%  "'
$ws.Range("E27").Formula = '="  -1.26%  "'
$ws.Range("D28").Formula = '="1.942.00"'
$ws.Range("E28").Formula = '="  -0.31%  "'
$ws.Range("D29").Formula = '="124.78"'
$ws.Range("D30").Formula = '="2.026"'
$ws.Range("E30").Formula = '="  -4.76%  "'
$ws.Range("E31").Formula = '="  -5.27%  "'
$ws.Range("D32").Formula = '="0.09048"'
$ws.Range("E32").Formula = '="  +1.79%  "'
$ws.Range("D33").Formula = '="3.644"'
$ws.Range("E33").Formula = '="  +0.04%  "'
$ws.Range("D34").Formula = '="5.368"'
$ws.Range("E34").Formula = '="  -2.89%  "'
$ws.Range("E35").Formula = '="  -0.61%  "'
$ws.Range("D36").Formula = '="11.61"'
$ws.Range("E36").Formula = '="  -4.81%  "'
$ws.Range("D37").Formula = '="0.05988"'
$ws.Range("E37").Formula = '="  -0.10%  "'
$ws.Range("E38").Formula = '="  -1.77%  "'
$ws.Range("D39").Formula = '="0.6232"'
$ws.Range("E39").Formula = '="  -1.57%  "'
$ws.Range("D40").Formula = '="4.868"'
$ws.Range("E40").Formula = '="  -1.59%  "'
$ws.Range("D41").Formula = '="1.181"'
$ws.Range("E41").Formula = '="  -0.20%  "'
$ws.Range("D42").Formula = '="1.375"'
$ws.Range("E42").Formula = '="  -0.73%  "'
$ws.Range("E43").Formula = '="  -2.39%  "'
$ws.Range("D44").Formula = '="13.10"'
$ws.Range("E44").Formula = '="  -2.09%  "'
$ws.Range("D45").Formula = '="3.701"'
$ws.Range("E45").Formula = '="  +0.28%  "'
$ws.Range("D46").Formula = '="0.5785"'
$ws.Range("E46").Formula = '="  -1.31%  "'
$ws.Range("D47").Formula = '="121.85"'
$ws.Range("E47").Formula = '="  -0.32%  "'
$ws.Range("D48").Formula = '="1.922"'
$ws.Range("E48").Formula = '="  -2.40%  "'
$ws.Range("D49").Formula = '="0.06838"'
$ws.Range("E49").Formula = '="  +0.52%  "'
$ws.Range("E50").Formula = '="  -4.88%  "'
$ws.Range("D51").Formula = '="71.13"'
$ws.Range("E51").Formula = '="  -2.58%  "'

$usedRange = $ws.Range("A1:E51")
$usedRange.Copy()
$usedRange.PasteSpecial(-4163)

